$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.964.34'
$ws.Range("E2").Value = '  +5.09%  '
$ws.Range("D3").Value = '3.506.09'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.89'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +7.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.505.36'
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.575'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.29'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").Value = '  +5.54%  '
$ws.Range("E12").Value = '  +3.94%  '
$ws.Range("D13").Value = '4.112.49'
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.134'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.29'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.43%  '
$ws.Range("E16").Value = '  +4.47%  '
$ws.Range("D17").Value = '66.922.79'
$ws.Range("E17").Value = '  +4.96%  '
$ws.Range("D18").Value = '3.488.84'
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("E19").Value = '  +3.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.05'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '394.54'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.77%  '
$ws.Range("E22").Value = '  +2.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000127'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +11.60%  '
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("E27").Value = '  +4.22%  '
$ws.Range("E28").Value = '  +2.39%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.37'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.36%  '
$ws.Range("E31").Value = '  +5.83%  '
$ws.Range("E32").Value = '  +4.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.47'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.49%  '
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  +6.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.77'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.900'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.36%  '
$ws.Range("E39").Value = '  +6.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0753'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.68'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.47%  '
$ws.Range("D42").Value = '2.847.48'
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.44'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.16%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.50'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.52'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0315'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("E48").Value = '  +7.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '349.56'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.08'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.85'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +12.79%  '
